# Auto-generated script applying scheduled price-data refresh to Siren_Profits workbook sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1330.7142
$ws.Range("J88").Value = 1447.6666
$ws.Range("L88").Value = 1447.6666
$ws.Range("N88").Value = -2259.6666

$ws.Range("H91").Value = 1330.7142
$ws.Range("J91").Value = 1447.6666
$ws.Range("L91").Value = 1447.6666
$ws.Range("N91").Value = -4255.6666

$ws.Range("H116").Value = 13892638
$ws.Range("I116").Value = 37039870
$ws.Range("J116").Value = 4297.8
$ws.Range("K116").Value = 37039870
$ws.Range("L116").Value = 4297.8
$ws.Range("M116").Value = -37036428
$ws.Range("N116").Value = -11181.8

$ws.Range("H132").Value = 3707996.2
$ws.Range("I132").Value = 4396.5835
$ws.Range("J132").Value = 33336794
$ws.Range("K132").Value = 13189.7505
$ws.Range("L132").Value = 100010382
$ws.Range("M132").Value = -10659.7505
$ws.Range("N132").Value = -100015442

$ws.Range("H137").Value = 8870.807000000001
$ws.Range("I137").Value = 10829.131
$ws.Range("J137").Value = 3240.625
$ws.Range("K137").Value = 32487.393
$ws.Range("L137").Value = 9721.875
$ws.Range("M137").Value = -29937.393
$ws.Range("N137").Value = -14821.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3590.057
$ws.Range("I32").Value = 3359.723
$ws.Range("K32").Value = 3359.723
$ws.Range("M32").Value = -3072.723

$ws.Range("H45").Value = 344313.5
$ws.Range("J45").Value = 11333.333
$ws.Range("L45").Value = 11333.333
$ws.Range("N45").Value = -12087.333

$ws.Range("H61").Value = 11626.348
$ws.Range("I61").Value = 16216.583
$ws.Range("J61").Value = 6618.8184
$ws.Range("K61").Value = 16216.583
$ws.Range("L61").Value = 6618.8184
$ws.Range("M61").Value = -16004.583
$ws.Range("N61").Value = -7042.8184

$ws.Range("H74").Value = 6659
$ws.Range("I74").Value = 11964.777
$ws.Range("K74").Value = 11964.777
$ws.Range("M74").Value = -11090.777

$ws.Range("H77").Value = 6659
$ws.Range("I77").Value = 11964.777
$ws.Range("K77").Value = 59823.885
$ws.Range("M77").Value = -55455.885

$ws.Range("H80").Value = 57249.5
$ws.Range("J80").Value = 84499
$ws.Range("L80").Value = 84499
$ws.Range("N80").Value = -86495

$ws.Range("H83").Value = 57249.5
$ws.Range("J83").Value = 84499
$ws.Range("L83").Value = 253497
$ws.Range("N83").Value = -263481

$ws.Range("H122").Value = 1207612
$ws.Range("I122").Value = 7450.5
$ws.Range("J122").Value = 2735090.2
$ws.Range("K122").Value = 22351.5
$ws.Range("L122").Value = 8205270.600000001
$ws.Range("M122").Value = -19901.5
$ws.Range("N122").Value = -8210170.600000001

$ws.Range("H132").Value = 2384.02
$ws.Range("I132").Value = 1345.5
$ws.Range("K132").Value = 4036.5
$ws.Range("M132").Value = -1506.5

$ws.Range("H136").Value = 11626.348
$ws.Range("I136").Value = 16216.583
$ws.Range("J136").Value = 6618.8184
$ws.Range("K136").Value = 48649.749
$ws.Range("L136").Value = 19856.4552
$ws.Range("M136").Value = -46099.749
$ws.Range("N136").Value = -24956.4552

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5141.8184
$ws.Range("I20").Value = 3127.1667
$ws.Range("J20").Value = 7559.4
$ws.Range("K20").Value = 3127.1667
$ws.Range("L20").Value = 7559.4
$ws.Range("M20").Value = -2880.1667
$ws.Range("N20").Value = -8053.4

$ws.Range("H140").Value = 94998.664
$ws.Range("J140").Value = 94998.664
$ws.Range("L140").Value = 94998.664
$ws.Range("N140").Value = -105358.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6331.9644
$ws.Range("I31").Value = 6512.2915
$ws.Range("K31").Value = 6512.2915
$ws.Range("M31").Value = -6217.2915

$ws.Range("H34").Value = 6331.9644
$ws.Range("I34").Value = 6512.2915
$ws.Range("K34").Value = 6512.2915
$ws.Range("M34").Value = -6310.2915

$ws.Range("H96").Value = 78020310
$ws.Range("J96").Value = 78020310
$ws.Range("L96").Value = 78020310
$ws.Range("N96").Value = -78025802

$ws.Range("H105").Value = 144463.36
$ws.Range("I105").Value = 155499
$ws.Range("K105").Value = 155499
$ws.Range("M105").Value = -153752

$ws.Range("H122").Value = 7521.5
$ws.Range("I122").Value = 7521.5
$ws.Range("K122").Value = 22564.5
$ws.Range("M122").Value = -20114.5

$ws.Range("H134").Value = 9942.166999999999
$ws.Range("I134").Value = 12561.889
$ws.Range("J134").Value = 2083
$ws.Range("K134").Value = 37685.667
$ws.Range("L134").Value = 6249
$ws.Range("M134").Value = -35150.667
$ws.Range("N134").Value = -11319

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 20711.7
$ws.Range("J22").Value = 50446.75
$ws.Range("L22").Value = 151340.25
$ws.Range("N22").Value = -151678.25

$ws.Range("H27").Value = 20711.7
$ws.Range("J27").Value = 50446.75
$ws.Range("L27").Value = 151340.25
$ws.Range("N27").Value = -151544.25

$ws.Range("H113").Value = 2657.5
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = ""

$ws.Range("H131").Value = 1932.2528
$ws.Range("I131").Value = 1164.8
$ws.Range("J131").Value = 2027
$ws.Range("K131").Value = 3494.4
$ws.Range("L131").Value = 6081
$ws.Range("M131").Value = 1545.6
$ws.Range("N131").Value = -16161

$ws.Range("H132").Value = 51894.7
$ws.Range("I132").Value = 1189.4
$ws.Range("J132").Value = 102600
$ws.Range("K132").Value = 10704.6
$ws.Range("L132").Value = 923400
$ws.Range("M132").Value = -8174.6
$ws.Range("N132").Value = -928460

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 70000
$ws.Range("I63").Value = 90000
$ws.Range("J63").Value = 50000
$ws.Range("K63").Value = 90000
$ws.Range("L63").Value = 50000
$ws.Range("M63").Value = -89314
$ws.Range("N63").Value = -51372

$ws.Range("H66").Value = 70000
$ws.Range("I66").Value = 90000
$ws.Range("J66").Value = 50000
$ws.Range("K66").Value = 270000
$ws.Range("L66").Value = 150000
$ws.Range("M66").Value = -266568
$ws.Range("N66").Value = -156864

$ws.Range("H80").Value = 5207.7827
$ws.Range("I80").Value = 5972.9287
$ws.Range("J80").Value = 4017.5557
$ws.Range("K80").Value = 5972.9287
$ws.Range("L80").Value = 4017.5557
$ws.Range("M80").Value = -4974.9287
$ws.Range("N80").Value = -6013.5557

$ws.Range("H83").Value = 5207.7827
$ws.Range("I83").Value = 5972.9287
$ws.Range("J83").Value = 4017.5557
$ws.Range("K83").Value = 29864.6435
$ws.Range("L83").Value = 20087.7785
$ws.Range("M83").Value = -24872.6435
$ws.Range("N83").Value = -30071.7785

$ws.Range("H113").Value = 17733.285
$ws.Range("I113").Value = 28283.25
$ws.Range("K113").Value = 28283.25
$ws.Range("M113").Value = -26113.25

$ws.Range("H122").Value = 9327.115
$ws.Range("I122").Value = 6267.778
$ws.Range("K122").Value = 18803.334
$ws.Range("M122").Value = -16353.334

$ws.Range("H132").Value = 4056.513
$ws.Range("I132").Value = 4133
$ws.Range("K132").Value = 12399
$ws.Range("M132").Value = -9869

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 10109
$ws.Range("I22").Value = 17166.5
$ws.Range("J22").Value = 1640
$ws.Range("K22").Value = 17166.5
$ws.Range("L22").Value = 1640
$ws.Range("M22").Value = -16871.5
$ws.Range("N22").Value = -2230

$ws.Range("H27").Value = 10109
$ws.Range("I27").Value = 17166.5
$ws.Range("J27").Value = 1640
$ws.Range("K27").Value = 17166.5
$ws.Range("L27").Value = 1640
$ws.Range("M27").Value = -17059.5
$ws.Range("N27").Value = -1854

$ws.Range("H56").Value = 11036
$ws.Range("I56").Value = 11036
$ws.Range("K56").Value = 11036
$ws.Range("M56").Value = -10345

$ws.Range("H122").Value = 6026.6553
$ws.Range("I122").Value = 5567.5264
$ws.Range("K122").Value = 16702.5792
$ws.Range("M122").Value = -14252.5792

$ws.Range("H136").Value = 5641.9473
$ws.Range("I136").Value = 2410
$ws.Range("K136").Value = 7230
$ws.Range("M136").Value = -4680

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7770.0625
$ws.Range("I81").Value = 11028.3
$ws.Range("J81").Value = 2339.6667
$ws.Range("K81").Value = 22056.6
$ws.Range("L81").Value = 4679.3334
$ws.Range("M81").Value = -20995.6
$ws.Range("N81").Value = -6801.3334

$ws.Range("H84").Value = 7770.0625
$ws.Range("I84").Value = 11028.3
$ws.Range("J84").Value = 2339.6667
$ws.Range("K84").Value = 110283
$ws.Range("L84").Value = 23396.667
$ws.Range("M84").Value = -104979
$ws.Range("N84").Value = -34004.667

$ws.Range("H122").Value = 4175.1133
$ws.Range("I122").Value = 2725.4688
$ws.Range("J122").Value = 6384.095
$ws.Range("K122").Value = 8176.4064
$ws.Range("L122").Value = 19152.285
$ws.Range("M122").Value = -5726.4064
$ws.Range("N122").Value = -24052.285

$ws.Range("H126").Value = 14641.277
$ws.Range("I126").Value = 18320.111
$ws.Range("J126").Value = 3604.7778
$ws.Range("K126").Value = 54960.333
$ws.Range("L126").Value = 10814.3334
$ws.Range("M126").Value = -52490.333
$ws.Range("N126").Value = -15754.3334

$ws.Range("H132").Value = 3031.3547
$ws.Range("I132").Value = 2294.0222
$ws.Range("J132").Value = 4983.1177
$ws.Range("K132").Value = 6882.0666
$ws.Range("L132").Value = 14949.3531
$ws.Range("M132").Value = -4352.0666
$ws.Range("N132").Value = -20009.3531

$ws.Range("H136").Value = 653503.25
$ws.Range("I136").Value = 708196.75
$ws.Range("J136").Value = 51875
$ws.Range("K136").Value = 2124590.25
$ws.Range("L136").Value = 155625
$ws.Range("M136").Value = -2122040.25
$ws.Range("N136").Value = -160725
